$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Copy the formatting of two "template" rows down onto the new rows so the
#    new cells inherit the same cell styles used throughout the sheet:
#      - row 569  -> G column has text (style s="1" t="s")
#      - row 563  -> G column is blank (style s="2")
# ---------------------------------------------------------------------------
$ws.Range("A569:I569").Copy()
$ws.Range("A571:I574").PasteSpecial(-4122)
$ws.Range("A576:I577").PasteSpecial(-4122)
$ws.Range("A580:I581").PasteSpecial(-4122)
$ws.Range("A583:I585").PasteSpecial(-4122)

$ws.Range("A563:I563").Copy()
$ws.Range("A570:I570").PasteSpecial(-4122)
$ws.Range("A575:I575").PasteSpecial(-4122)
$ws.Range("A578:I579").PasteSpecial(-4122)
$ws.Range("A582:I582").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Fill in the new training-log rows (570-585).
# ---------------------------------------------------------------------------

# Row 570
$ws.Range("A570").Value = 45979
$ws.Range("B570").Value = "Maé Clavel"
$ws.Range("C570").Value = 60
$ws.Range("D570").Value = 2
$ws.Range("E570").Value = 1
$ws.Range("F570").Value = 0
$ws.Range("H570").Value = 10
$ws.Range("I570").Formula = "=C570*D570"

# Row 571
$ws.Range("A571").Value = 45979
$ws.Range("B571").Value = "Yoan Zouma"
$ws.Range("C571").Value = 60
$ws.Range("D571").Value = 2
$ws.Range("E571").Value = 5
$ws.Range("F571").Value = 5
$ws.Range("G571").Value = "Ischio"
$ws.Range("H571").Value = 5
$ws.Range("I571").Formula = "=C571*D571"

# Row 572
$ws.Range("A572").Value = 45979
$ws.Range("B572").Value = "Yoann Martelat"
$ws.Range("C572").Value = 60
$ws.Range("D572").Value = 3
$ws.Range("E572").Value = 3
$ws.Range("F572").Value = 3
$ws.Range("G572").Value = "Genou"
$ws.Range("H572").Value = 6
$ws.Range("I572").Formula = "=C572*D572"

# Row 573
$ws.Range("A573").Value = 45979
$ws.Range("B573").Value = "Kamal Bafounta"
$ws.Range("C573").Value = 60
$ws.Range("D573").Value = 3
$ws.Range("E573").Value = 1
$ws.Range("F573").Value = 2
$ws.Range("G573").Value = "Genou"
$ws.Range("H573").Value = 3
$ws.Range("I573").Formula = "=C573*D573"

# Row 574
$ws.Range("A574").Value = 45979
$ws.Range("B574").Value = "Omar Benyounes"
$ws.Range("C574").Value = 60
$ws.Range("D574").Value = 5
$ws.Range("E574").Value = 6
$ws.Range("F574").Value = 1
$ws.Range("G574").Value = "Ischio droit"
$ws.Range("H574").Value = 6
$ws.Range("I574").Formula = "=C574*D574"

# Row 575
$ws.Range("A575").Value = 45979
$ws.Range("B575").Value = "Ilyes Boughanmi"
$ws.Range("C575").Value = 60
$ws.Range("D575").Value = 6
$ws.Range("E575").Value = 6
$ws.Range("F575").Value = 0
$ws.Range("H575").Value = 0
$ws.Range("I575").Formula = "=C575*D575"

# Row 576
$ws.Range("A576").Value = 45979
$ws.Range("B576").Value = "Naim Ighbane"
$ws.Range("C576").Value = 60
$ws.Range("D576").Value = 3
$ws.Range("E576").Value = 4
$ws.Range("F576").Value = 3
$ws.Range("G576").Value = "Genou"
$ws.Range("H576").Value = 3
$ws.Range("I576").Formula = "=C576*D576"

# Row 577
$ws.Range("A577").Value = 45979
$ws.Range("B577").Value = "Karim Belmahi"
$ws.Range("C577").Value = 60
$ws.Range("D577").Value = 5
$ws.Range("E577").Value = 5
$ws.Range("F577").Value = 3
$ws.Range("G577").Value = "Mollet "
$ws.Range("H577").Value = 10
$ws.Range("I577").Formula = "=C577*D577"

# Row 578
$ws.Range("A578").Value = 45979
$ws.Range("B578").Value = "Jeremie Laurent"
$ws.Range("C578").Value = 60
$ws.Range("D578").Value = 5
$ws.Range("E578").Value = 5
$ws.Range("F578").Value = 0
$ws.Range("H578").Value = 7
$ws.Range("I578").Formula = "=C578*D578"

# Row 579
$ws.Range("A579").Value = 45979
$ws.Range("B579").Value = "Emmanuel Valey"
$ws.Range("C579").Value = 60
$ws.Range("D579").Value = 4
$ws.Range("E579").Value = 2
$ws.Range("F579").Value = 0
$ws.Range("H579").Value = 9
$ws.Range("I579").Formula = "=C579*D579"

# Row 580
$ws.Range("A580").Value = 45979
$ws.Range("B580").Value = "Karahali Souaré"
$ws.Range("C580").Value = 60
$ws.Range("D580").Value = 5
$ws.Range("E580").Value = 5
$ws.Range("F580").Value = 6
$ws.Range("G580").Value = "Cheville"
$ws.Range("H580").Value = 5
$ws.Range("I580").Formula = "=C580*D580"

# Row 581
$ws.Range("A581").Value = 45979
$ws.Range("B581").Value = "Sofiane Belle"
$ws.Range("C581").Value = 60
$ws.Range("D581").Value = 5
$ws.Range("E581").Value = 7
$ws.Range("F581").Value = 3
$ws.Range("G581").Value = "Ischio"
$ws.Range("H581").Value = 5
$ws.Range("I581").Formula = "=C581*D581"

# Row 582
$ws.Range("A582").Value = 45979
$ws.Range("B582").Value = "Mattheo Haon"
$ws.Range("C582").Value = 60
$ws.Range("D582").Value = 6
$ws.Range("E582").Value = 3
$ws.Range("F582").Value = 0
$ws.Range("H582").Value = 8
$ws.Range("I582").Formula = "=C582*D582"

# Row 583
$ws.Range("A583").Value = 45979
$ws.Range("B583").Value = "Levy Ndoutoume"
$ws.Range("C583").Value = 60
$ws.Range("D583").Value = 5
$ws.Range("E583").Value = 6
$ws.Range("F583").Value = 6
$ws.Range("G583").Value = "Béquille "
$ws.Range("H583").Value = 8
$ws.Range("I583").Formula = "=C583*D583"

# Row 584
$ws.Range("A584").Value = 45979
$ws.Range("B584").Value = "Wael Fareh"
$ws.Range("C584").Value = 60
$ws.Range("D584").Value = 4
$ws.Range("E584").Value = 3
$ws.Range("F584").Value = 4
$ws.Range("G584").Value = "Genou"
$ws.Range("H584").Value = 8
$ws.Range("I584").Formula = "=C584*D584"

# Row 585
$ws.Range("A585").Value = 45979
$ws.Range("B585").Value = "Hedi Nasri"
$ws.Range("C585").Value = 60
$ws.Range("D585").Value = 4
$ws.Range("E585").Value = 2
$ws.Range("F585").Value = 1
$ws.Range("G585").Value = "Hanche"
$ws.Range("H585").Value = 9
$ws.Range("I585").Formula = "=C585*D585"

# ---------------------------------------------------------------------------
# 3. Update the sheet view: scroll position and active selection.
# ---------------------------------------------------------------------------
$aw = $excel.ActiveWindow
$aw.ScrollRow = 554
$aw.ScrollColumn = 1
$ws.Range("L565").Select()

Write-Output "done"
